# Add two new data rows (Elva Lynn / Daphne Boone) to the group_b sheet,
# replacing the two stray, empty formatted rows (16 and 21) that previously
# only carried a one-off fill style with no actual data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old rows 16 and 21 only held direct (unused) cell formatting - no
# values. Clear() wipes both content and formatting so the rows disappear
# from sheetData entirely (mirrors them being dropped from the saved file).
$ws.Range("A16:G16").Clear() | Out-Null
$ws.Range("A21:G21").Clear() | Out-Null

# New data row 2: Elva Lynn
$ws.Range("A2").Value = "Elva"
$ws.Range("B2").Value = "Lynn"
$ws.Range("C2").Value = "Elva.Lynn1987@gmail.com"
$ws.Range("D2").Value = "C9nqK4Xz6gbUocjf"
$ws.Range("E2").Value = "81.28.96.172:43739"
$ws.Range("F2").Value = "a4ZJCKXpyPxqtgAt"
$ws.Range("G2").Value = "WZRYymEeVjfQeR6Z"

# New data row 3: Daphne Boone (shares proxy user/pass with row 2)
$ws.Range("A3").Value = "Daphne"
$ws.Range("B3").Value = "Boone"
$ws.Range("C3").Value = "Daphne.Boone2002@gmail.com"
$ws.Range("D3").Value = "aPxbySqRts82761Y"
$ws.Range("E3").Value = "81.28.96.172:30660"
$ws.Range("F3").Value = "a4ZJCKXpyPxqtgAt"
$ws.Range("G3").Value = "WZRYymEeVjfQeR6Z"

# Match the saved selection/active cell recorded in the workbook.
$ws.Range("L20").Select() | Out-Null
